$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Add the new Operator import to the "Import" value cell (B2)
$importCell = $ws.Range("B2")
$importCell.Value = $importCell.Value() + ",org.openmrs.module.drools.calculation.Operator"

# 2. Update the CONDITION formulas to use calculationService.checkMostRecentObs(...) with the Operator enum
$ws.Range("B8").Value = "calculationService.checkMostRecentObs(`$patient, SYSTOLIC_UUID, Operator.GTE, `$param)"
$ws.Range("C8").Value = "calculationService.checkMostRecentObs(`$patient, SYSTOLIC_UUID, Operator.LT, `$param)"
$ws.Range("D8").Value = "calculationService.checkMostRecentObs(`$patient, DIASTOLIC_UUID, Operator.GTE, `$param)"
$ws.Range("E8").Value = "calculationService.checkMostRecentObs(`$patient, DIASTOLIC_UUID, Operator.LT, `$param)"

# 3. Normalize the duplicated cell style used by the blank "separator" row (row 4) and the
#    threshold value cells (rows 10-12) so it collapses onto the same underlying style as the
#    rest of the numeric-right-aligned cells (mirrors Excel's own style-table de-duplication).
$ws.Range("D1").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$ws.Range("B4").PasteSpecial(-4122)
$ws.Range("B10:E12").PasteSpecial(-4122)
